$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing time portion from the date strings in column A
$ws.Range("A1").Value = "15/03/2019"
$ws.Range("A2").Value = "22/03/2019"
$ws.Range("A3").Value = "23/04/2019"

# Update the active selection
$ws.Range("A12").Select()
